# Update Leave Card 1/15/2024 3:59 PM
#
# This script applies the 2024 leave-card entries to "ARCILLA, MAYETTE A."
# leave card workbook:
#   - Fills in VL/SL earned amounts for Oct/Nov/Dec 2023 rows (67-69)
#   - Adds a 3 day absence and "FL(3-0-0)" particular for Dec 2023 (row 69)
#   - Inserts a new "2024" year-header row (row 70)
#   - Populates the monthly period dates for Jan-Dec 2024 (rows 71-82)
#   - Records a "SP(1-0-0)" particular and "PARENTAL 01/19/2024" remark for Jan 2024 (row 71)
#   - Extends the leave table by one blank trailing row (128 -> 129)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# ---------------------------------------------------------------------------
# 1. Extend Table1 / sheetData with a new trailing blank row (129), pushing
#    the heavier "final row" border formatting down from row 128 to row 129,
#    and giving row 128 the regular interior-row formatting.
# ---------------------------------------------------------------------------
$ws.Range("A128:K128").Copy($ws.Range("A129:K129"))
$lo.Resize($ws.Range("A8:K129"))
$ws.Range("G129").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$ws.Range("A127:K127").Copy()
$ws.Range("A128:K128").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill in VL/SL earned credits for the existing Oct/Nov/Dec 2023 rows.
# ---------------------------------------------------------------------------
# Row 67 (10/1/2023): EARNED 1.25
$ws.Range("C67").Value = 1.25

# Row 68 (11/1/2023): EARNED 1.25
$ws.Range("C68").Value = 1.25

# Row 69 (12/1/2023): EARNED 1.25, Absence W/Pay 3 days
$ws.Range("C69").Value = 1.25
$ws.Range("D69").Value = 3

# ---------------------------------------------------------------------------
# 3. Insert the "2024" year marker row (row 70), matching the styling used by
#    the other year-header rows (bold, quote-prefixed text in the date column).
# ---------------------------------------------------------------------------
$ws.Range("A70").Value = "'2024"
$ws.Range("A70").Font.Bold = $true

# Row 69 particular (added after the "2024" marker so the shared string
# table ends up appended in the same order as the authored edit).
$ws.Range("B69").Value = "FL(3-0-0)"

# ---------------------------------------------------------------------------
# 4. Populate the monthly period rows for 2024 (rows 71-82).
# ---------------------------------------------------------------------------
$ws.Range("A71").Value = 45292
$ws.Range("B71").Value = "SP(1-0-0)"
$ws.Range("K71").Value = "PARENTAL 01/19/2024"

$ws.Range("A72").Value = 45323
$ws.Range("A73").Value = 45352
$ws.Range("A74").Value = 45383
$ws.Range("A75").Value = 45413
$ws.Range("A76").Value = 45444
$ws.Range("A77").Value = 45474
$ws.Range("A78").Value = 45505
$ws.Range("A79").Value = 45536
$ws.Range("A80").Value = 45566
$ws.Range("A81").Value = 45597
$ws.Range("A82").Value = 45627

# ---------------------------------------------------------------------------
# 5. Restore the selection to reflect the last-edited cell.
# ---------------------------------------------------------------------------
$ws.Range("G75").Select()

$wb.Save()
